# Apply the "ger ita jpn prophet" update:
#  - add two new columns: "Targeted Month for Prediction" (E) and "Source model" (F)
#  - fill in previously-empty Germany/Italy/Japan prediction rows
#  - populate the new E/F columns for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): new column headers -------------------------------
$ws.Range("E2").Value = "Targeted Month for Prediction"
$ws.Range("F2").Value = "Source model"

# --- Row 3: Canada ----------------------------------------------------------
$ws.Range("E3").Value = "May"
$ws.Range("F3").Value = "Prophet CV"

# --- Row 4: France (values were stored as text before; make them numeric) --
$ws.Range("B4").Value = 6.991768
$ws.Range("C4").Value = 7.241449
$ws.Range("D4").Value = 6.729585
$ws.Range("E4").Value = "April"
$ws.Range("F4").Value = "Prophet CV"

# --- Row 5: Germany (previously blank) --------------------------------------
$ws.Range("B5").Value = 3.1088488
$ws.Range("C5").Value = 3.344162
$ws.Range("D5").Value = 2.851482
$ws.Range("E5").Value = "May"
$ws.Range("F5").Value = "Prophet CV"

# --- Row 6: Italy (previously blank) ----------------------------------------
$ws.Range("B6").Value = 8.14796
$ws.Range("C6").Value = 8.513808
$ws.Range("D6").Value = 7.817619
$ws.Range("E6").Value = "April"
$ws.Range("F6").Value = "Prophet CV"

# --- Row 7: Japan (previously blank) ----------------------------------------
$ws.Range("B7").Value = 2.008909
$ws.Range("C7").Value = 2.293228
$ws.Range("D7").Value = 2.013057
$ws.Range("E7").Value = "April"
$ws.Range("F7").Value = "Prophet CV"

# --- Row 8: United Kingdom ---------------------------------------------------
$ws.Range("E8").Value = "March"
$ws.Range("F8").Value = "Prophet CV"

# --- Row 9: United States -----------------------------------------------------
$ws.Range("E9").Value = "May"
$ws.Range("F9").Value = "Prophet CV"

# --- Row heights: every row goes from the old default (12.8) to 15 ---------
for ($r = 1; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# --- Selection moves to B8 per the saved view -------------------------------
$ws.Range("B8").Select() | Out-Null
